$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 583.2105
$ws.Range("I12").Value = 424.1
$ws.Range("J12").Value = 760
$ws.Range("K12").Value = 424.1
$ws.Range("L12").Value = 760
$ws.Range("M12").Value = -254.1
$ws.Range("N12").Value = -1100
$ws.Range("H33").Value = 142.9
$ws.Range("I33").Value = 146.07143
$ws.Range("K33").Value = 146.07143
$ws.Range("M33").Value = 82.92857000000001
$ws.Range("H62").Value = 2665.8333
$ws.Range("I62").Value = 2799
$ws.Range("K62").Value = 2799
$ws.Range("M62").Value = -2175
$ws.Range("H65").Value = 2665.8333
$ws.Range("I65").Value = 2799
$ws.Range("K65").Value = 13995
$ws.Range("M65").Value = -10875
$ws.Range("H86").Value = 1900.3334
$ws.Range("I86").Value = 1900.75
$ws.Range("K86").Value = 1900.75
$ws.Range("M86").Value = -777.75
$ws.Range("H89").Value = 1900.3334
$ws.Range("I89").Value = 1900.75
$ws.Range("K89").Value = 9503.75
$ws.Range("M89").Value = -3887.75
$ws.Range("H92").Value = 336
$ws.Range("I92").Value = 288.7143
$ws.Range("K92").Value = 288.7143
$ws.Range("M92").Value = 959.2857
$ws.Range("H125").Value = 2258.3333
$ws.Range("J125").Value = 2758.75
$ws.Range("L125").Value = 24828.75
$ws.Range("N125").Value = -29748.75
$ws.Range("H135").Value = 923
$ws.Range("I135").Value = 894.9
$ws.Range("K135").Value = 8054.099999999999
$ws.Range("M135").Value = -5519.099999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 397919.28
$ws.Range("I2").Value = 794553.3
$ws.Range("K2").Value = 794553.3
$ws.Range("M2").Value = -794440.3
$ws.Range("H32").Value = 3636.82
$ws.Range("I32").Value = 3552.3435
$ws.Range("K32").Value = 3552.3435
$ws.Range("M32").Value = -3265.3435
$ws.Range("H45").Value = 10001590
$ws.Range("I45").Value = 45000610
$ws.Range("K45").Value = 45000610
$ws.Range("M45").Value = -45000233
$ws.Range("H61").Value = 40421.19
$ws.Range("I61").Value = 57699.785
$ws.Range("K61").Value = 57699.785
$ws.Range("M61").Value = -57487.785
$ws.Range("H80").Value = 49333.332
$ws.Range("J80").Value = 49333.332
$ws.Range("L80").Value = 49333.332
$ws.Range("N80").Value = -51329.332
$ws.Range("H83").Value = 49333.332
$ws.Range("J83").Value = 49333.332
$ws.Range("L83").Value = 147999.996
$ws.Range("N83").Value = -157983.996
$ws.Range("H97").Value = 1102.9375
$ws.Range("I97").Value = 1101.9231
$ws.Range("J97").Value = 1107.3334
$ws.Range("K97").Value = 1101.9231
$ws.Range("L97").Value = 1107.3334
$ws.Range("M97").Value = -605.9231
$ws.Range("N97").Value = -2099.3334
$ws.Range("H110").Value = 937.17645
$ws.Range("I110").Value = 708.1539
$ws.Range("K110").Value = 708.1539
$ws.Range("M110").Value = 1336.8461
$ws.Range("H116").Value = 397919.28
$ws.Range("I116").Value = 794553.3
$ws.Range("K116").Value = 794553.3
$ws.Range("M116").Value = -792259.3
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("H136").Value = 40421.19
$ws.Range("I136").Value = 57699.785
$ws.Range("K136").Value = 173099.355
$ws.Range("M136").Value = -170549.355
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 397919.28
$ws.Range("I3").Value = 794553.3
$ws.Range("K3").Value = 794553.3
$ws.Range("M3").Value = -794439.3
$ws.Range("H20").Value = 4598.2
$ws.Range("I20").Value = 2997
$ws.Range("K20").Value = 2997
$ws.Range("M20").Value = -2750
$ws.Range("H22").Value = 711.2222
$ws.Range("I22").Value = 583.5
$ws.Range("K22").Value = 583.5
$ws.Range("M22").Value = -410.5
$ws.Range("H86").Value = 1000875.5
$ws.Range("I86").Value = 1501
$ws.Range("K86").Value = 1501
$ws.Range("M86").Value = -378
$ws.Range("H89").Value = 1000875.5
$ws.Range("I89").Value = 1501
$ws.Range("K89").Value = 7505
$ws.Range("M89").Value = -1889
$ws.Range("H105").Value = 2135.8
$ws.Range("I105").Value = 2089.8057
$ws.Range("K105").Value = 2089.8057
$ws.Range("M105").Value = -342.8056999999999
$ws.Range("H134").Value = 4753.2812
$ws.Range("I134").Value = 4672.4585
$ws.Range("J134").Value = 4995.75
$ws.Range("K134").Value = 14017.3755
$ws.Range("L134").Value = 14987.25
$ws.Range("M134").Value = -11482.3755
$ws.Range("N134").Value = -20057.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1417.6666
$ws.Range("I16").Value = 1417.6666
$ws.Range("K16").Value = 1417.6666
$ws.Range("M16").Value = -1130.6666
$ws.Range("H22").Value = 1205.5333
$ws.Range("I22").Value = 430
$ws.Range("K22").Value = 430
$ws.Range("M22").Value = -80
$ws.Range("H31").Value = 2301.9697
$ws.Range("I31").Value = 2183.889
$ws.Range("J31").Value = 2443.6667
$ws.Range("K31").Value = 2183.889
$ws.Range("L31").Value = 2443.6667
$ws.Range("M31").Value = -1888.889
$ws.Range("N31").Value = -3033.6667
$ws.Range("H34").Value = 2301.9697
$ws.Range("I34").Value = 2183.889
$ws.Range("J34").Value = 2443.6667
$ws.Range("K34").Value = 2183.889
$ws.Range("L34").Value = 2443.6667
$ws.Range("M34").Value = -1981.889
$ws.Range("N34").Value = -2847.6667
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = $null
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H105").Value = 836.7143
$ws.Range("I105").Value = 921.5
$ws.Range("K105").Value = 921.5
$ws.Range("M105").Value = 825.5
$ws.Range("H107").Value = 4902
$ws.Range("I107").Value = 1382.3334
$ws.Range("K107").Value = 1382.3334
$ws.Range("M107").Value = 537.6666
$ws.Range("H113").Value = 1417.6666
$ws.Range("I113").Value = 1417.6666
$ws.Range("K113").Value = 1417.6666
$ws.Range("M113").Value = 752.3334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13912467
$ws.Range("J131").Value = 27317.193
$ws.Range("L131").Value = 81951.579
$ws.Range("N131").Value = -92031.579
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1587.1666
$ws.Range("I102").Value = 1686.1818
$ws.Range("J102").Value = 498
$ws.Range("K102").Value = 1686.1818
$ws.Range("L102").Value = 498
$ws.Range("M102").Value = -64.18180000000007
$ws.Range("N102").Value = -3742
$ws.Range("H113").Value = 1383.1666
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1499.8
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 1499.8
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -5839.8
$ws.Range("H126").Value = 3271187.8
$ws.Range("I126").Value = 3706890.5
$ws.Range("J126").Value = 3416.5
$ws.Range("K126").Value = 11120671.5
$ws.Range("L126").Value = 10249.5
$ws.Range("M126").Value = -11118201.5
$ws.Range("N126").Value = -15189.5
$ws.Range("H136").Value = 19877.777
$ws.Range("J136").Value = 19877.777
$ws.Range("L136").Value = 59633.33099999999
$ws.Range("N136").Value = -64733.33099999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 21333.334
$ws.Range("J50").Value = 27000
$ws.Range("L50").Value = 27000
$ws.Range("N50").Value = -28274
$ws.Range("H88").Value = 18000
$ws.Range("I88").Value = 18000
$ws.Range("K88").Value = 18000
$ws.Range("M88").Value = -17572
$ws.Range("H91").Value = 18000
$ws.Range("I91").Value = 18000
$ws.Range("K91").Value = 18000
$ws.Range("M91").Value = -16518
$ws.Range("H93").Value = 1961.6666
$ws.Range("I93").Value = 1046.75
$ws.Range("J93").Value = 3791.5
$ws.Range("K93").Value = 1046.75
$ws.Range("L93").Value = 3791.5
$ws.Range("M93").Value = 201.25
$ws.Range("N93").Value = -6287.5
$ws.Range("H122").Value = 2901.4546
$ws.Range("I122").Value = 2734
$ws.Range("K122").Value = 8202
$ws.Range("M122").Value = -5752
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47092
$ws.Range("J70").Value = 47092
$ws.Range("L70").Value = 47092
$ws.Range("N70").Value = -47722
$ws.Range("H73").Value = 47092
$ws.Range("J73").Value = 47092
$ws.Range("L73").Value = 47092
$ws.Range("N73").Value = -49276
$ws.Range("H139").Value = 69833.336
$ws.Range("J139").Value = 69833.336
$ws.Range("L139").Value = 69833.336
$ws.Range("N139").Value = -80113.336
